$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "59.811.95") that Excel would
# otherwise auto-convert to a number on assignment. Force text storage by
# pre-formatting the column as Text, then strip the formatting back off so
# the cells keep the workbook's original (unstyled) appearance.
$dCol = $ws.Range("D2:D51")
$dCol.NumberFormat = "@"

$ws.Range("D2").Value = "59.811.95"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.373.20"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "555.77"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "133.51"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "0.105"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "5.65"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").Value = "24.39"
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("D14").Value = "2.801.18"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "59.755.99"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "2.373.13"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "11.10"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "4.47"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "321.13"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "64.12"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "0.173"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "8.39"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "1.80"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "0.0₃0757"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "170.08"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "6.06"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +9.49%  "
$ws.Range("D33").Value = "0.398"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").Value = "18.15"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "1.32"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "317.91"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").Value = "38.60"
$ws.Range("D42").Value = "144.88"
$ws.Range("E42").Value = "  +4.19%  "
$ws.Range("E43").Value = "  -3.83%  "
$ws.Range("D44").Value = "0.0969"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "19.62"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "0.0509"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").Value = "0.570"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "0.0217"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").Value = "11.05"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "4.67"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -1.59%  "

$dCol.ClearFormats()
